$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginPage")

# Update the expected-message cell text. The shared string
# "User-ID must not be blank" is no longer used anywhere, so it will
# naturally drop out of the workbook's shared string table once this is
# the only reference rewritten.
$ws.Range("C2").Value = "User or Password is not valid"

# Move / record the active selection to C2 (previously F12).
$ws.Activate()
$ws.Range("C2").Select()
